$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 (existing CNPJ rows 5-161 shift down to 6-162).
$ws.Rows("5:5").Insert()

# New row 5 holds a CNPJ-validator seed/sample value.
$c5 = $ws.Range("A5")
$c5.Value = 123
$ws.Rows("5:5").RowHeight = 16.5

# Format the new cell: 14-digit zero-padded numeric mask, left aligned,
# no border, default font (matches the sheet's normal/blank-row style).
$c5.Style = "Normal"
$c5.NumberFormat = "00000000000000"
$c5.HorizontalAlignment = -4131

Write-Host "done"
